$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) and sheet "全部类型" (sheet4) share the same
# underlying event records; both need their "想去人数" (F) / "最低票价" (G)
# numbers refreshed to the newly scraped values.

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("G2").Value = 50

    $ws.Range("F4").Value = 37
    $ws.Range("F5").Value = 3406
    $ws.Range("F6").Value = 2150
    $ws.Range("F7").Value = 412

    if ($name -eq "展览") {
        $ws.Range("F9").Value = 48
        $ws.Range("F10").Value = 35
        $ws.Range("F11").Value = 1261
        $ws.Range("F13").Value = 1534
        $ws.Range("F14").Value = 111
    } else {
        $ws.Range("F10").Value = 48
        $ws.Range("F11").Value = 35
        $ws.Range("F14").Value = 1261
        $ws.Range("F16").Value = 1534
        $ws.Range("F17").Value = 111
    }
}
